# Timetracking workbook update: add "adminservice tests" timings and two new
# "frontend" tasks (admin verify/block/unblock, parameter edit) under the
# AdminPanel user story.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50 ("adminservice tests"): fill in estimated / actual time.
$ws.Range("C50").Value = "15min"
$ws.Range("D50").Value = "10min"

# Row 51 currently holds the lone "frontend" placeholder in column B.
# Free it first so the shared string slot it occupies can be repurposed for
# the new "1h 15min" duration value (set into D51), matching how the
# original authors re-used that string table slot.
$ws.Range("B51").Value = ""
$ws.Range("D51").Value = "1h 15min"

# Now give row 51 its real task name and estimated time.
$ws.Range("B51").Value = "frontend - admin verify, block, unblock"
$ws.Range("C51").Value = "1h"

# Add the new row 52 for the second frontend task.
$ws.Range("B52").Value = "frontend - parameter edit"

# Update the view so the newly added row is visible/selected, matching the
# author's saved view state as closely as the object model allows.
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
[void]$ws.Range("B52").Select()
